# Auto-generated edit script: refresh market-price derived columns (H-N)
# in the leve profit tracking sheets, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5 - Met a Sticky End
$ws.Range("H5").Value = 234
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents() | Out-Null
# Row 11 - Gotta Bounce
$ws.Range("H11").Value = 300
$ws.Range("I11").Value = 300
$ws.Range("K11").Value = 300
$ws.Range("M11").Value = -160
# Row 18 - You Grow, Girl
$ws.Range("H18").Value = 871
$ws.Range("I18").Value = 871
$ws.Range("K18").Value = 871
$ws.Range("M18").Value = -587
# Row 21 - Book and a Hard Place
$ws.Range("H21").Value = 25951
$ws.Range("I21").Value = 29131
$ws.Range("K21").Value = 29131
$ws.Range("M21").Value = -28663
# Row 23 - There's Something about Bury
$ws.Range("H23").Value = 25951
$ws.Range("I23").Value = 29131
$ws.Range("K23").Value = 29131
$ws.Range("M23").Value = -28897
# Row 32 - Automata for the People
$ws.Range("H32").Value = 744.7857
$ws.Range("I32").Value = 656.2
$ws.Range("J32").Value = 794
$ws.Range("K32").Value = 656.2
$ws.Range("L32").Value = 794
$ws.Range("M32").Value = -330.2
$ws.Range("N32").Value = -1446
# Row 40 - Stuck in the Moment
$ws.Range("H40").Value = 1967.1428
$ws.Range("I40").Value = 1835.7142
$ws.Range("J40").Value = 2230
$ws.Range("K40").Value = 1835.7142
$ws.Range("L40").Value = 2230
$ws.Range("M40").Value = -1660.7142
$ws.Range("N40").Value = -2580
# Row 43 - Growing Is Knowing
$ws.Range("H43").Value = 1918.1538
$ws.Range("I43").Value = 1932.8572
$ws.Range("J43").Value = 1901
$ws.Range("K43").Value = 1932.8572
$ws.Range("L43").Value = 1901
$ws.Range("M43").Value = -1863.8572
$ws.Range("N43").Value = -2039
# Row 51 - A Bile Business
$ws.Range("H51").Value = 2814.6667
$ws.Range("I51").Value = 1834
$ws.Range("J51").Value = 2937.25
$ws.Range("K51").Value = 1834
$ws.Range("L51").Value = 2937.25
$ws.Range("M51").Value = -1350
$ws.Range("N51").Value = -3905.25
# Row 55 - A Real Smooth Move
$ws.Range("H55").Value = 295.66666
$ws.Range("I55").Value = 295.66666
$ws.Range("K55").Value = 295.66666
$ws.Range("M55").Value = -81.66665999999998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5 - The Alloyed Truth
$ws.Range("H5").Value = 916.6667
$ws.Range("I5").Value = 916.6667
$ws.Range("K5").Value = 916.6667
$ws.Range("M5").Value = -804.6667
# Row 32 - Ingot We Trust
$ws.Range("H32").Value = 18186064
$ws.Range("I32").Value = 20412296
$ws.Range("J32").Value = 5166.6665
$ws.Range("K32").Value = 20412296
$ws.Range("L32").Value = 5166.6665
$ws.Range("M32").Value = -20412009
$ws.Range("N32").Value = -5740.6665

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4 - Mending Fences
$ws.Range("H4").Value = 916.6667
$ws.Range("I4").Value = 916.6667
$ws.Range("K4").Value = 916.6667
$ws.Range("M4").Value = -801.6667
# Row 22 - Riveting Run
$ws.Range("I22").Value = 550
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = -377
$ws.Range("N22").Value = -446

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 8 - Bows for the Boys
$ws.Range("H8").Value = 1046
$ws.Range("J8").Value = 1046
$ws.Range("L8").Value = 1046
$ws.Range("N8").Value = -1326
# Row 19 - Shielding Sales
$ws.Range("H19").Value = 466.66666
$ws.Range("I19").Value = 300
$ws.Range("K19").Value = 300
$ws.Range("M19").Value = -130
# Row 24 - What You Need
$ws.Range("H24").Value = 466.66666
$ws.Range("I24").Value = 300
$ws.Range("K24").Value = 300
$ws.Range("M24").Value = -130
# Row 25 - Bowing to Necessity
$ws.Range("H25").Value = 70000
$ws.Range("J25").Value = 70000
$ws.Range("L25").Value = 70000
$ws.Range("N25").Value = -70348
# Row 51 - Greenstone for Greenhorns
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 30000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 30000
$ws.Range("M51").ClearContents() | Out-Null
$ws.Range("N51").Value = -31472
# Row 56 - Trident and Error
$ws.Range("H56").Value = 5359.3
$ws.Range("I56").Value = 3796.5
$ws.Range("J56").Value = 5750
$ws.Range("K56").Value = 3796.5
$ws.Range("L56").Value = 5750
$ws.Range("M56").Value = -2951.5
$ws.Range("N56").Value = -7440
# Row 60 - Bowing to Greater Power
$ws.Range("H60").Value = 10735.792
$ws.Range("I60").Value = 2749.5
$ws.Range("J60").Value = 11461.818
$ws.Range("K60").Value = 2749.5
$ws.Range("L60").Value = 11461.818
$ws.Range("M60").Value = -2238.5
$ws.Range("N60").Value = -12483.818
# Row 61 - Incant Now, Think Later
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 30000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 30000
$ws.Range("M61").ClearContents() | Out-Null
$ws.Range("N61").Value = -30696

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5 - What a Sap
$ws.Range("H5").Value = 659.61536
$ws.Range("I5").Value = 396.25
$ws.Range("K5").Value = 1188.75
$ws.Range("M5").Value = -1076.75
# Row 122 - Salt of the North
$ws.Range("H122").Value = 942.1111
$ws.Range("I122").Value = 370
$ws.Range("J122").Value = 1399.8
$ws.Range("K122").Value = 3330
$ws.Range("L122").Value = 12598.2
$ws.Range("M122").Value = -880
$ws.Range("N122").Value = -17498.2
# Row 135 - Not-so-secret Ingredient
$ws.Range("H135").Value = 659.61536
$ws.Range("I135").Value = 396.25
$ws.Range("K135").Value = 3566.25
$ws.Range("M135").Value = -1031.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 9 - From the Sands to the Stage
$ws.Range("H9").Value = 652.5
$ws.Range("I9").Value = 299.5
$ws.Range("J9").Value = 1005.5
$ws.Range("K9").Value = 299.5
$ws.Range("L9").Value = 1005.5
$ws.Range("M9").Value = -75.5
$ws.Range("N9").Value = -1453.5
# Row 22 - Skin off Their Backs
$ws.Range("H22").Value = 2459.375
$ws.Range("I22").Value = 1382.1428
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 1382.1428
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -1087.1428
$ws.Range("N22").Value = -10590
# Row 27 - Fire and Hide
$ws.Range("H27").Value = 2459.375
$ws.Range("I27").Value = 1382.1428
$ws.Range("J27").Value = 10000
$ws.Range("K27").Value = 1382.1428
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = -1275.1428
$ws.Range("N27").Value = -10214
# Row 46 - Supply Side Logic
$ws.Range("H46").Value = 1896.1111
$ws.Range("I46").Value = 2317.4
$ws.Range("J46").Value = 1369.5
$ws.Range("K46").Value = 2317.4
$ws.Range("L46").Value = 1369.5
$ws.Range("M46").Value = -2129.4
$ws.Range("N46").Value = -1745.5
# Row 100 - Tiger in the Sack
$ws.Range("H100").Value = 1215
$ws.Range("I100").Value = 1150.8334
$ws.Range("J100").Value = 1600
$ws.Range("K100").Value = 1150.8334
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = -609.8334
$ws.Range("N100").Value = -2682

